{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The first paragraph in the document is \"Each method dumps its data in a\n// CSV file.\" - insert the two new bullet items immediately before it, so\n// they become the new first two paragraphs. insertParagraph() on an\n// existing list paragraph inherits its pPr (ListParagraph style / numPr),\n// matching the target list formatting.\nconst firstParagraph = paragraphs.items[0];\n\nfirstParagraph.insertParagraph(\n  \"Added the standard Copyright header \",\n  Word.InsertLocation.before\n);\nfirstParagraph.insertParagraph(\n  \"License: Apache 2.0\",\n  Word.InsertLocation.before\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The first paragraph currently in the document is \"Each method dumps its\n# data in a CSV file.\" - insert the two new bullet items before it, each\n# inheriting that paragraph's list style (ListParagraph / numId 1).\n$firstPara = $d.Paragraphs(1)\n$firstPara.Range.InsertParagraphBefore()\n$d.Paragraphs(1).Range.Text = \"Added the standard Copyright header \"\n\n$secondPara = $d.Paragraphs(2)\n$secondPara.Range.InsertParagraphBefore()\n$d.Paragraphs(2).Range.Text = \"License: Apache 2.0\"\n"}
